$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Transmitance values (column B, rows 3-18) from 1 to 100
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 2).Value = 100
}

# Update the selected cell to B18
$ws.Range("B18").Select()
